$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 288
$ws.Cells.Item(3, 6).Value = 656
$ws.Cells.Item(6, 6).Value = 2850
$ws.Cells.Item(10, 6).Value = 591
$ws.Cells.Item(14, 6).Value = 6037
$ws.Cells.Item(15, 6).Value = 650
$ws.Cells.Item(16, 6).Value = 1057
$ws.Cells.Item(18, 6).Value = 247
$ws.Cells.Item(21, 6).Value = 575
$ws.Cells.Item(22, 6).Value = 11
$ws.Cells.Item(23, 6).Value = 49
$ws.Cells.Item(24, 6).Value = 32
$ws.Cells.Item(25, 6).Value = 140
$ws.Cells.Item(26, 6).Value = 1334
$ws.Cells.Item(28, 6).Value = 1012
$ws.Cells.Item(29, 6).Value = 60
$ws.Cells.Item(30, 6).Value = 2078
$ws.Cells.Item(31, 6).Value = 192
$ws.Cells.Item(32, 6).Value = 20
$ws.Cells.Item(35, 6).Value = 3343

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(3, 6).Value = 27
$ws.Cells.Item(7, 6).Value = 98
$ws.Cells.Item(11, 6).Value = 651
$ws.Cells.Item(19, 6).Value = 65
$ws.Cells.Item(21, 6).Value = 357
$ws.Cells.Item(27, 6).Value = 158
$ws.Cells.Item(29, 6).Value = 72

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(4, 6).Value = 44
$ws.Cells.Item(5, 6).Value = 2604
$ws.Cells.Item(6, 6).Value = 1151
$ws.Cells.Item(8, 6).Value = 1504
$ws.Cells.Item(9, 6).Value = 420
$ws.Cells.Item(10, 6).Value = 119
$ws.Cells.Item(11, 6).Value = 5
$ws.Cells.Item(12, 6).Value = 669

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(4, 6).Value = 2604
$ws.Cells.Item(5, 6).Value = 1151
$ws.Cells.Item(6, 6).Value = 1504
$ws.Cells.Item(7, 6).Value = 420
$ws.Cells.Item(8, 6).Value = 119
$ws.Cells.Item(9, 6).Value = 288
$ws.Cells.Item(10, 6).Value = 656
$ws.Cells.Item(11, 6).Value = 2850
$ws.Cells.Item(13, 6).Value = 669
$ws.Cells.Item(14, 6).Value = 591
$ws.Cells.Item(15, 6).Value = 98
$ws.Cells.Item(19, 6).Value = 6038
$ws.Cells.Item(21, 6).Value = 650
$ws.Cells.Item(22, 6).Value = 1057
$ws.Cells.Item(23, 6).Value = 247
$ws.Cells.Item(26, 6).Value = 575
$ws.Cells.Item(30, 6).Value = 65
$ws.Cells.Item(31, 6).Value = 11
$ws.Cells.Item(32, 6).Value = 32
$ws.Cells.Item(34, 6).Value = 357
$ws.Cells.Item(38, 6).Value = 158
$ws.Cells.Item(40, 6).Value = 60
$ws.Cells.Item(41, 6).Value = 72
$ws.Cells.Item(43, 6).Value = 2078
$ws.Cells.Item(46, 6).Value = 192
$ws.Cells.Item(47, 6).Value = 20
$ws.Cells.Item(49, 6).Value = 3343
